$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 20

# Column A holds a date-looking value but is stored as plain text in this
# workbook (matching the existing rows). Temporarily force a text number
# format so Excel doesn't auto-convert the string into a date serial
# number, then restore the default "Normal" style so no stray per-cell
# style is left behind.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12/14/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 12298.94
$ws.Cells.Item($row, 3).Value = 0.2055609289118182
$ws.Cells.Item($row, 4).Value = 0.7944390710881818
$ws.Cells.Item($row, 5).Value = -132.4
$ws.Cells.Item($row, 6).Value = -28.78
$ws.Cells.Item($row, 7).Value = -20746
$ws.Cells.Item($row, 8).Value = -67.98
$ws.Cells.Item($row, 9).Value = -424.67
$ws.Cells.Item($row, 10).Value = -14.38
